$d = $word.ActiveDocument

$d.Content.Find.Execute("80×59=4720", $true, $false, $false, $false, $false, $true, 1, $false, "50×14=700", 2) | Out-Null
$d.Content.Find.Execute("14×47=658", $true, $false, $false, $false, $false, $true, 1, $false, "59×55=3245", 2) | Out-Null
$d.Content.Find.Execute("48×17=816", $true, $false, $false, $false, $false, $true, 1, $false, "50×58=2900", 2) | Out-Null
$d.Content.Find.Execute("81×53=4293", $true, $false, $false, $false, $false, $true, 1, $false, "84×62=5208", 2) | Out-Null
$d.Content.Find.Execute("50×11=550", $true, $false, $false, $false, $false, $true, 1, $false, "55×11=605", 2) | Out-Null
$d.Content.Find.Execute("53×88=4664", $true, $false, $false, $false, $false, $true, 1, $false, "86×14=1204", 2) | Out-Null
$d.Content.Find.Execute("54×97=5238", $true, $false, $false, $false, $false, $true, 1, $false, "85×91=7735", 2) | Out-Null
$d.Content.Find.Execute("48×90=4320", $true, $false, $false, $false, $false, $true, 1, $false, "73×29=2117", 2) | Out-Null
$d.Content.Find.Execute("53×84=4452", $true, $false, $false, $false, $false, $true, 1, $false, "68×38=2584", 2) | Out-Null
$d.Content.Find.Execute("22×79=1738", $true, $false, $false, $false, $false, $true, 1, $false, "98×96=9408", 2) | Out-Null
$d.Content.Find.Execute("96×81=7776", $true, $false, $false, $false, $false, $true, 1, $false, "90×72=6480", 2) | Out-Null
$d.Content.Find.Execute("85×92=7820", $true, $false, $false, $false, $false, $true, 1, $false, "97×29=2813", 2) | Out-Null
$d.Content.Find.Execute("20×17=340", $true, $false, $false, $false, $false, $true, 1, $false, "46×93=4278", 2) | Out-Null
$d.Content.Find.Execute("86×96=8256", $true, $false, $false, $false, $false, $true, 1, $false, "50×63=3150", 2) | Out-Null
$d.Content.Find.Execute("60×28=1680", $true, $false, $false, $false, $false, $true, 1, $false, "82×48=3936", 2) | Out-Null
$d.Content.Find.Execute("79×91=7189", $true, $false, $false, $false, $false, $true, 1, $false, "84×41=3444", 2) | Out-Null
$d.Content.Find.Execute("94×91=8554", $true, $false, $false, $false, $false, $true, 1, $false, "28×45=1260", 2) | Out-Null
$d.Content.Find.Execute("96×83=7968", $true, $false, $false, $false, $false, $true, 1, $false, "64×92=5888", 2) | Out-Null
$d.Content.Find.Execute("13×19=247", $true, $false, $false, $false, $false, $true, 1, $false, "18×50=900", 2) | Out-Null
$d.Content.Find.Execute("36×53=1908", $true, $false, $false, $false, $false, $true, 1, $false, "82×97=7954", 2) | Out-Null
$d.Content.Find.Execute("98×13=1274", $true, $false, $false, $false, $false, $true, 1, $false, "29×14=406", 2) | Out-Null
$d.Content.Find.Execute("72×30=2160", $true, $false, $false, $false, $false, $true, 1, $false, "88×44=3872", 2) | Out-Null
$d.Content.Find.Execute("79×37=2923", $true, $false, $false, $false, $false, $true, 1, $false, "71×19=1349", 2) | Out-Null
$d.Content.Find.Execute("48×98=4704", $true, $false, $false, $false, $false, $true, 1, $false, "66×34=2244", 2) | Out-Null
$d.Content.Find.Execute("41×35=1435", $true, $false, $false, $false, $false, $true, 1, $false, "76×86=6536", 2) | Out-Null
